# Adds a new "2021" column (column R) to the statistics table on the
# active sheet, mirroring the formatting of the preceding "2020" column
# (column Q), fills in the new year's data, and moves the active
# selection to the top of the newly added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at R, shifting cells right and carrying the
# formatting from the column immediately to its left (column Q) so the
# new column's cells inherit the same styles (borders, fonts, number
# formats) as the existing year columns instead of ending up blank.
$ws.Columns.Item(18).Insert(-4121)

# Populate the new column's header (year) and data value.
$ws.Range("R4").Value2 = 2021
$ws.Range("R5").Value2 = 72

# Move the selection to the new column's top cell, matching where the
# user would naturally land after adding the new column of data.
$ws.Range("R1").Select()
